# NeetCode 150 Tracker — "new Trapping rain water added"
# Fills in Week-2 rows (rows 8-14) with dates, progress notes & status,
# widens column H, and updates the saved sheet view (zoom + selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NeetCode 150 Tracker")

function Set-DateCell($addr, [int]$year, [int]$month, [int]$day) {
    # Clone the number format of an existing date cell (A2 already carries
    # the workbook's date style) so we reuse the same cellXf instead of
    # minting a brand-new one, then stamp the date value itself.
    $ws.Range("A2").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $ws.Range($addr).Value = (Get-Date -Year $year -Month $month -Day $day).Date
}

function Set-NoteCell($addr, $text) {
    # Clone the wrap-text style already used by the Notes column (H6) so
    # new notes line up with the existing cellXf reuse pattern.
    $ws.Range("H6").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $ws.Range($addr).Value = $text
}

# ---- Row 8 : Encode and Decode Strings ----
Set-DateCell "A8" 2026 1 13
$ws.Range("G8").Value = " 🟡"
Set-NoteCell "H8" " main logic here is length-prefixed serialization"
$ws.Rows.Item(8).RowHeight = 16

# ---- Row 9 : Longest Consecutive Sequence ----
Set-DateCell "A9" 2026 1 14
$ws.Range("F9").Value = "day8"
$ws.Range("G9").Value = "❌"
Set-NoteCell "H9" "do it later"
$ws.Rows.Item(9).RowHeight = 16

# ---- Row 10 : Valid Palindrome ----
Set-DateCell "A10" 2026 1 14
$ws.Range("F10").Value = "day 9"
$ws.Range("G10").Value = " 🟡"
Set-NoteCell "H10" "done only important logic is regex and ascii to convert it/[a-z0-9]/i.test(s[r])"
$ws.Rows.Item(10).RowHeight = 32

# ---- Row 11 : Two Sum II ----
Set-DateCell "A11" 2026 1 15
$ws.Range("F11").Value = "day9"
$ws.Range("G11").Value = "✅ 🟡"
Set-NoteCell "H11" "solve self by bruteForce but now go for optimal"
$ws.Rows.Item(11).RowHeight = 16

# ---- Row 12 : 3Sum ----
Set-DateCell "A12" 2026 1 15
$ws.Range("F12").Value = "day9"
$ws.Range("G12").Value = "❌"

# ---- Row 13 : Container With Most Water ----
Set-DateCell "A13" 2026 1 25
$ws.Range("F13").Value = "day10"
$ws.Range("G13").Value = "✅"
Set-NoteCell "H13" "solve by slef because two pointer clear "
$ws.Rows.Item(13).RowHeight = 16

# ---- Row 14 : Trapping Rain Water ----
Set-DateCell "A14" 2026 1 25
$ws.Range("F14").Value = "day10"

# Row 5 note cell shrank from 64pt to 48pt tall now that other rows share
# the load of the Week-1 notes.
$ws.Rows.Item(5).RowHeight = 48

# Column H (Notes) widened to fit the new commentary.
$ws.Columns.Item(8).ColumnWidth = 48.33

# Saved view: bumped zoom and moved the active selection to the new
# Trapping Rain Water row.
$ws.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 165
$ws.Range("G15").Select() | Out-Null
